$wb = $excel.ActiveWorkbook

# --- Work on the "ok (2)" sheet (soon to be renamed "CURATED") ---
$ws3 = $wb.Worksheets.Item("ok (2)")
$ws3.Select()

# New header + formula column C, built in the same three chunks Excel itself
# would have used (one single-cell entry, then two fill/enter blocks), so the
# shared-formula grouping in the saved XML matches what a real edit produces.
$ws3.Range("C1").Value = "string1"
$ws3.Range("C2").Formula = "=CONCATENATE(""['"",A2,""', "",INT(B2),""],"")"
$ws3.Range("C3:C66").Formula = "=CONCATENATE(""['"",A3,""', "",INT(B3),""],"")"
$ws3.Range("C67:C70").Formula = "=CONCATENATE(""['"",A67,""', "",INT(B67),""],"")"

# Sort all three columns descending by the CO2 value column (B), header row excluded.
$sortRange = $ws3.Range("A1:C70")
$sortRange.Sort($ws3.Range("B1"), 2)

# Turn on AutoFilter across the curated table.
$ws3.Range("A1:C70").AutoFilter()

# AutoFilter creates a hidden, sheet-scoped _FilterDatabase defined name.
$filterName = $ws3.Names.Add("_xlnm._FilterDatabase", "=CURATED!`$A`$1:`$C`$70")
$filterName.Visible = $false

# Leave the curated selection on column C, matching the last manual selection.
$ws3.Range("C2:C70").Select()

# Rename the sheet now that it holds the curated/sorted data.
$ws3.Name = "CURATED"

# --- Restore focus to Sheet1, which becomes the active tab again ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Select()
$ws1.Range("M16").Select()
